$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 'signatures' column (header in J1, no data below it) is being
# dropped entirely. Deleting the whole column shifts everything that
# was to its right (the 'expansion' column, previously K) one place
# to the left so it becomes the new column J - exactly what the
# target diff shows.
$ws.Columns("J:J").Delete()

# Excel leaves the selection on the cell that slid into the deleted
# column's place.
$ws.Range("J1").Select()
